$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- GPS STATUS (1Hz) data ID 5002 section ---
# Previously rows 23-24 were a merged "Vertical dilution of precision" bit-field
# (1 bit value + 7 bit decimeters exponent). The commit replaces that with a new
# 2-bit "GPS advanced fix" field (RTK status) plus a 6-bit RESERVED field.

# Unmerge the old B23:B24 merged label cell
$ws.Range("B23:B24").UnMerge()

# Row 23: GPS advanced fix
$ws.Range("B23").Value = "GPS advanced fix"
$ws.Range("C23").Value = "N/A"
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = "0: no advanced fix, 1: DGPS, 2: RTK_FLOAT, 3: RTK_FIXED"

# Row 24: RESERVED bits
$ws.Range("B24").Value = "RESERVED"
$ws.Range("D24").Value = 6
$ws.Range("C24").Clear()
$ws.Range("E24").Clear()

Write-Host "Done"
